{"js": "// REPORTGEN-665: fix typos in security reports\n//\n// This applies the three deliberate text edits made by the author:\n//   1. \"STIG standards\"  -> \"CISQ security standards\"\n//   2. \"This section provide a summary\" -> \"This section provides a summary\"\n//   3. \"mesurement\"      -> \"measurement\"\n//\n// (The surrounding bookmark-id renumbering / VML shapetype re-serialisation /\n// latentStyles churn visible in the source diff are artefacts Word's engine\n// regenerates automatically on every save and are not reachable through the\n// Word JavaScript API's content object model, so they are intentionally left\n// alone here.)\n\nconst body = context.document.body;\n\n// 1. \"STIG standards\" -> \"CISQ security standards\"\nlet results = body.search(\"STIG standards\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"CISQ security standards\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. \"This section provide a summary\" -> \"This section provides a summary\"\nresults = body.search(\"This section provide a summary\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"This section provides a summary\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3. \"mesurement\" -> \"measurement\"\nresults = body.search(\"mesurement\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText(\"measurement\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# REPORTGEN-665: fix typos in security reports\n#\n# This applies the three deliberate text edits made by the author:\n#   1. \"STIG standards\"  -> \"CISQ security standards\"\n#   2. \"This section provide a summary\" -> \"This section provides a summary\"\n#   3. \"mesurement\"      -> \"measurement\"\n#\n# (The surrounding bookmark-id renumbering / VML shapetype re-serialisation /\n# latentStyles churn visible in the source diff are artefacts Word's engine\n# regenerates automatically on every save and are not exposed as settable\n# values anywhere in the Word object model, so they are intentionally left\n# alone here.)\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Text = $findText\n    $r.Find.MatchCase = $true\n    $r.Find.MatchWholeWord = $false\n    $r.Find.Replacement.Text = $replaceText\n    # wdFindContinue=1 (Wrap), wdReplaceAll=2 (Replace)\n    $r.Find.Execute($null, $true, $null, $null, $null, $null, $true, 1, $null, $replaceText, 2) | Out-Null\n}\n\nReplace-FirstMatch \"STIG standards\" \"CISQ security standards\"\nReplace-FirstMatch \"This section provide a summary\" \"This section provides a summary\"\nReplace-FirstMatch \"mesurement\" \"measurement\"\n"}
